$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/re-order header row with the new imaging columns
$ws.Range("A1").Value = "MRN"
$ws.Range("B1").Value = " Nombre"
$ws.Range("C1").Value = " Physician"
$ws.Range("D1").Value = " Appt Time"
$ws.Range("E1").Value = " WR Timestamp"
$ws.Range("F1").Value = " EX Timestamp"
$ws.Range("G1").Value = " FC Start"
$ws.Range("H1").Value = " FC End"
$ws.Range("I1").Value = "Imaging"
$ws.Range("J1").Value = "Imaging Timestamp"
$ws.Range("K1").Value = " DC Timestamp"
$ws.Range("L1").Value = " WR Total Time"
$ws.Range("M1").Value = " EX Total Time"
$ws.Range("N1").Value = "Total Time"
$ws.Range("O1").Value = " AT Entry"

# Column widths for the new/resized columns
$ws.Range("E:E").ColumnWidth = 14.7109375
$ws.Range("F:F").ColumnWidth = 13.85546875
$ws.Range("G:G").ColumnWidth = 13.28515625
$ws.Range("H:H").ColumnWidth = 12
$ws.Range("I:I").ColumnWidth = 9.42578125
$ws.Range("J:J").ColumnWidth = 18.42578125
$ws.Range("K:K").ColumnWidth = 14.140625
$ws.Range("L:L").ColumnWidth = 14.140625
$ws.Range("M:M").ColumnWidth = 13.28515625
$ws.Range("N:N").ColumnWidth = 12.7109375

# Selection matches the diff's sheetView selection
$ws.Range("D2").Select()
